# Update LR-pair TPM data: drop the "ECs" target-cluster row and refresh
# the remaining "Resolving-Mac" row with newly computed TPM-based values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 currently holds the "ECs" target-cluster data; row 3 holds the
# "Resolving-Mac" data. Delete row 2 (ECs) so Resolving-Mac becomes row 2,
# then overwrite its numeric columns with the recomputed TPM values.
$ws.Rows(2).Delete()

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06792633333333332
$ws.Range("H2").Value = 0.203779
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2121096666666666
$ws.Range("N2").Value = 0.6363289999999999
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.01440783192122222
$ws.Range("R2").Value = 0.129670487291
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
